$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; B=1.459612070389937;  C=3099.503889238888;  D=0.1575252929769615; E=8.660232485948974;  G=3109.781259088204}
    @{Row=3; B=3.230985683306322;  C=1.667794583268128;  D=0.8054896365839992; E=0.496779210170732;  G=6.201049113329182}
    @{Row=4; B=3.230985683306322;  C=1.667794583268128;  D=26.21740644021617;  E=0.496779210170732;  G=31.61296591696135}
    @{Row=5; B=3.230985683306322;  C=1.667794583268128;  D=0.8054896365839992; E=0.496779210170732;  G=6.201049113329182}
    @{Row=6; B=1.459612070389937;  C=1.667794583268128;  D=26.21740644021617;  E=0.496779210170732;  G=29.84159230404497}
    @{Row=7; B=0.3048080303191223; C=1.667794583268128;  D=3.900430680208489;  E=645.3272768299601;  G=651.2003101237558}
    @{Row=8; B=1.459612070389937;  C=1.667794583268128;  D=0.8054896365839992; E=0.496779210170732;  G=4.429675500412797}
    @{Row=9; B=3.230985683306322;  C=1.667794583268128;  D=0.8054896365839992; E=0.496779210170732;  G=6.201049113329182}
)

foreach ($r in $data) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}

$wb.Save()
